$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BOM rows (A:G, row 2-14) got reordered while researching part
# characteristics: the Kamaya part (old row 5) moved up above the Darfon
# part, and the AVX / Bourns parts (old rows 2-3) moved down to the
# bottom of the list. Re-write the table in its new row order.

$rows = @(
    @{ A = "FRC32C1J599JB";        B = "R1-R43";                   C = 84;  D = 0.31;               G = "Kamaya"  },
    @{ A = "C2012NP0460FFAQZ";     B = "C10-C121";                 C = 112; D = 0.2;                G = "Darfon"  },
    @{ A = "C1210C321K5PACTU";     B = "C205-C209";                C = 5;   D = 0.1;                G = "Kemet"   },
    @{ A = "RK73H2HAGTE2726D";     B = "R37-R136";                 C = 100; D = 0.1;                G = "KOA"     },
    @{ A = "GHOIJOIEFCNFJ32324";   B = "U1-U4";                    C = 5;   D = 0.56000000000000005; G = "UNKOWN"  },
    @{ A = "MA2225CG3R7A500";      B = "C210-C223";                C = 14;  D = 0.24;               G = "Meritek" },
    @{ A = "RK73H1JATTD9R90F";     B = "R137-R161";                C = 25;  D = 0.1;                G = "KOA"     },
    @{ A = "C0201C100M9PACTU";     B = "C224-C244";                C = 21;  D = 0.1;                G = "Kemet"   },
    @{ A = "RR0816P-721-F-M";      B = "R162-R269";                C = 108; D = 0.13;               G = "Susumu"  },
    @{ A = "C0603X5R1V225C160BB";  B = "C245";                     C = 1;   D = 0.14000000000000001; G = "TDK"     },
    @{ A = "RC0805JKP135K31L";     B = "R270- R355";               C = 86;  D = 0.1;                G = "YAGEO"   },
    @{ A = "KAM15GS71C563CU";      B = "C1, C2, C3, C4, C5-C9";    C = 9;   D = 0.21;               G = "AVX"     },
    @{ A = "CRM1206-JZ-3657ELF";   B = "R1-R36";                   C = 36;  D = 0.38;               G = "Bourns"  }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 7).Value = $row.G
    $r = $r + 1
}

# The view had scrolled down (topLeftCell=A4, selection B5); after the
# reordering work the user ended up with the view back at the top and a
# new selection further down the sheet.
$ws.Range("A1").Select()
$ws.Range("B20").Select()
